$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '26.751.61'
$c.ClearFormats()

$c = $ws.Range('E2')
$c.NumberFormat = "@"
$c.Value = '  +0.36%  '
$c.ClearFormats()

$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '1.603.12'
$c.ClearFormats()

$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  +0.39%  '
$c.ClearFormats()

$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  +0.21%  '
$c.ClearFormats()

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '211.89'
$c.ClearFormats()

$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  +0.23%  '
$c.ClearFormats()

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.515'
$c.ClearFormats()

$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  +0.11%  '
$c.ClearFormats()

$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  +0.19%  '
$c.ClearFormats()

$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  +0.12%  '
$c.ClearFormats()

$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  -0.07%  '
$c.ClearFormats()

$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  +0.53%  '
$c.ClearFormats()

$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  +0.70%  '
$c.ClearFormats()

$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '1.826.64'
$c.ClearFormats()

$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  +0.28%  '
$c.ClearFormats()

$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '1.605.34'
$c.ClearFormats()

$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  +0.17%  '
$c.ClearFormats()

$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  +0.69%  '
$c.ClearFormats()

$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  +0.04%  '
$c.ClearFormats()

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '65.07'
$c.ClearFormats()

$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  +0.08%  '
$c.ClearFormats()

$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.0₃0739'
$c.ClearFormats()

$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  +0.29%  '
$c.ClearFormats()

$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '209.93'
$c.ClearFormats()

$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  +0.42%  '
$c.ClearFormats()

$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  +0.22%  '
$c.ClearFormats()

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '7.15'
$c.ClearFormats()

$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  +0.24%  '
$c.ClearFormats()

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '2.24'
$c.ClearFormats()

$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  -3.86%  '
$c.ClearFormats()

$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  +0.04%  '
$c.ClearFormats()

$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '143.54'
$c.ClearFormats()

$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  -0.42%  '
$c.ClearFormats()

$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  +0.22%  '
$c.ClearFormats()

$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  -0.57%  '
$c.ClearFormats()

$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '15.34'
$c.ClearFormats()

$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  +0.43%  '
$c.ClearFormats()

$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  -1.23%  '
$c.ClearFormats()

$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  -0.02%  '
$c.ClearFormats()

$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '3.27'
$c.ClearFormats()

$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  +0.70%  '
$c.ClearFormats()

$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  +0.52%  '
$c.ClearFormats()

$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '1.291.88'
$c.ClearFormats()

$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  +0.30%  '
$c.ClearFormats()

$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  +0.56%  '
$c.ClearFormats()

$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  +0.41%  '
$c.ClearFormats()

$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.601'
$c.ClearFormats()

$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  -2.15%  '
$c.ClearFormats()

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '1.18'
$c.ClearFormats()

$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  +11.80%  '
$c.ClearFormats()

$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  +0.15%  '
$c.ClearFormats()

$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  -0.06%  '
$c.ClearFormats()

$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '5.42'
$c.ClearFormats()

$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  -1.31%  '
$c.ClearFormats()

$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  -0.22%  '
$c.ClearFormats()

$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  -0.02%  '
$c.ClearFormats()

$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '62.99'
$c.ClearFormats()

$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  -0.73%  '
$c.ClearFormats()

$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '1.738.55'
$c.ClearFormats()

$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  +0.26%  '
$c.ClearFormats()

$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  -0.34%  '
$c.ClearFormats()

$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  -1.01%  '
$c.ClearFormats()

$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  +0.30%  '
$c.ClearFormats()

$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  +1.45%  '
$c.ClearFormats()

$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  +0.18%  '
$c.ClearFormats()

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '7.47'
$c.ClearFormats()

$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  +1.33%  '
$c.ClearFormats()

$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  +1.00%  '
$c.ClearFormats()

